# Apply the "Edit new date and expect result" change to TestCasesRobustness.
#
# Semantics of this test sheet:
#   A = Spending, B = VisitsPerMonth, C = Points,
#   D = ExpectedRank, E = ActualRank, F = Result
#
# The commit bumps the VisitsPerMonth ("new date"/input) values used across
# the robustness test rows, and updates the expected-rank/result columns
# ("expect result") accordingly. The "Gold" rank tier no longer exists, so
# every expected "Gold" becomes "Standard", and row 20's scenario (which
# used to land on Standard/PASS) now lands on Silver/FAIL once A20 is
# corrected from 5000 to 50000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (VisitsPerMonth): 4 -> 15 across the bulk of the rows -------
$ws.Range("B2").Value  = 15
$ws.Range("B3").Value  = 15
$ws.Range("B4").Value  = 15
$ws.Range("B5").Value  = 15
$ws.Range("B6").Value  = 15
$ws.Range("B7").Value  = 15
$ws.Range("B8").Value  = 15

# B3:B8 also drop their wrap/vertical-center formatting, matching B2's plain
# style (the edit effectively propagated B2's formatting down).
$ws.Range("B3:B8").ClearFormats()

$ws.Range("B15").Value = 15
$ws.Range("B16").Value = 15
$ws.Range("B17").Value = 15
$ws.Range("B18").Value = 15
$ws.Range("B19").Value = 15
$ws.Range("B20").Value = 15

# --- Column B: the "Silver" boundary rows shift up too ---------------------
$ws.Range("B12").Value = 29
$ws.Range("B13").Value = 30
$ws.Range("B14").Value = 31

# --- Column D (ExpectedRank): "Gold" tier removed -> "Standard" ------------
$ws.Range("D5").Value  = "Standard"
$ws.Range("D6").Value  = "Standard"
$ws.Range("D7").Value  = "Standard"
$ws.Range("D18").Value = "Standard"
$ws.Range("D19").Value = "Standard"

# --- Row 20: fix the Spending input and its expected outcome ---------------
$ws.Range("A20").Value = 50000
$ws.Range("E20").Value = "Silver"
$ws.Range("F20").Value = "FAIL"

# --- UI state: last selected cell moved to J15 ------------------------------
$ws.Range("J15").Select()
